$d = $word.ActiveDocument
$newBodyXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Anotações</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>OneToMany</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>ManyToOne</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Exemplo Funcionario_Dependente: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>OneToMany: Um funcionário pode ter muitos dependentes? Sim</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>OneToMany: Um dependente pode ter muitos funcionários? Não</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>ManyToOne: Muitos Funcionários podem ter um dependente? Não</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>ManyToOne: Muitos Dependentes podem ter um funcionário? Sim</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Então na classe Funcionário utiliza o OneToMany com uma lista de dependentes</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Então na classe Dependente utiliza o ManyToOne com uma variável funcionário do tipo Funcionario</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>JoinColumn()</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Exemplo Funcionario_Dependente:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Utilizado para acrescentar a chave estrangeira, então deve ser usado no lado que receberá a chave estrangeira no caso o lado do “N”</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>Serializable</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Exemplos Funcionario_Dependente:</w:t>
      </w:r>
      <w:r>
        <w:br/>
        <w:t xml:space="preserve">A classe dependente possui uma chave composta baseada no cpf do funcionário “fcpf” e no nome do dependente, então na Classe @Embeddable (dependente_id) terá as duas chaves primarias</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Usado em classes @Embeddable quando se está utilizando chave primária COMPOSTA</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>MapsI</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
        <w:t>d</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:r>
        <w:t>Exemplo Funcionario_Dependente:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>De acordo com o @OneToMany e o @ManyToOne das duas classes, já vai saber quais classes estão acontecendo a conexão. O MapsId automaticamente pega a chave primária da outra classe e preenche a variável onde o nome está correspondente, exemplo:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>@MapsId(“funcionário_cpf”) -&gt; O “funcionário_cpf” está presente na Classe DependenteId que é uma classe que tem chave completa composta, então o maps automaticamente pega da Funcionario a chave primaria e joga os dados dentro da variável “funcionário_cpf” da classe DependenteId</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:bCs/>
        </w:rPr>
      </w:pPr>
    </w:p>
'@
$d.Content.InsertXML($newBodyXml)
